# Update cryptos worksheet with the latest ranking snapshot values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.838.21"
$ws.Range("D3").Value = "1.887.43"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'0.7525"
$ws.Range("E5").Value = "  -3.82%  "
$ws.Range("D6").Value = "'242.33"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("D9").Value = "'25.34"
$ws.Range("E9").Value = "  -1.58%  "
$ws.Range("D10").Value = "'0.07123"
$ws.Range("E10").Value = "  -3.05%  "
$ws.Range("D11").Value = "'0.08483"
$ws.Range("E11").Value = "  +4.83%  "
$ws.Range("D12").Value = "'0.7603"
$ws.Range("E12").Value = "  -1.76%  "
$ws.Range("D13").Value = "1.896.36"
$ws.Range("E13").Value = "  -1.48%  "
$ws.Range("E14").Value = "  -2.84%  "
$ws.Range("D16").Value = "'6.138"
$ws.Range("E16").Value = "  -1.89%  "
$ws.Range("D17").Value = "29.801.45"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").Value = "'13.72"
$ws.Range("E18").Value = "  -1.94%  "
$ws.Range("D19").Value = "'243.80"
$ws.Range("E19").Value = "  -1.42%  "
$ws.Range("D20").Value = "'0.000007815"
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").Value = "'0.9994"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.134.68"
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("D23").Value = "'8.000"
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Value = "'0.1598"
$ws.Range("E25").Value = "  +0.16%  "
$ws.Range("D26").Value = "'9.373"
$ws.Range("E26").Value = "  -0.91%  "
$ws.Range("D27").Value = "'162.42"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").Value = "'2.030"
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("E30").Value = "  +4.00%  "
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("D32").Value = "'4.514"
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("D33").Value = "'4.127"
$ws.Range("E33").Value = "  +1.50%  "
$ws.Range("D34").Value = "'0.05425"
$ws.Range("E34").Value = "  -2.71%  "
$ws.Range("D35").Value = "'1.242"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Value = "'0.7505"
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("D37").Value = "'1.002"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").Value = "'2.711"
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("D39").Value = "'0.01949"
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("D40").Value = "'2.775"
$ws.Range("E40").Value = "  -0.97%  "
$ws.Range("D41").Value = "'0.4457"
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("E42").Value = "  +2.45%  "
$ws.Range("D43").Value = "1.090.37"
$ws.Range("E43").Value = "  -2.37%  "
$ws.Range("D44").Value = "'72.65"
$ws.Range("E44").Value = "  -2.51%  "
$ws.Range("D45").Value = "'0.8587"
$ws.Range("E45").Value = "  +0.87%  "
$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "'7.751"
$ws.Range("E47").Value = "  +2.96%  "
$ws.Range("D48").Value = "'102.43"
$ws.Range("E48").Value = "  -0.15%  "
$ws.Range("D49").Value = "'1.860"
$ws.Range("E49").Value = "  -1.65%  "
$ws.Range("D50").Value = "'3.062"
$ws.Range("E50").Value = "  +0.31%  "
$ws.Range("D51").Value = "2.042.68"
$ws.Range("E51").Value = "  +0.37%  "
